$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.455.93"
$ws.Range("E2").Value = "  +4.96%  "
$ws.Range("D3").Value = "1.723.61"
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").Value = "'225.66"
$ws.Range("E5").Value = "  +3.23%  "
$ws.Range("D6").Value = "'0.5365"
$ws.Range("E6").Value = "  +2.54%  "
$ws.Range("D7").Value = "'1.005"
$ws.Range("E7").Value = "  +0.12%  "
$ws.Range("D8").Value = "'0.2679"
$ws.Range("E8").Value = "  +0.52%  "
$ws.Range("D9").Value = "'0.06589"
$ws.Range("E9").Value = "  +3.89%  "
$ws.Range("D10").Value = "'21.62"
$ws.Range("E10").Value = "  +5.32%  "
$ws.Range("D11").Value = "'0.07763"
$ws.Range("E11").Value = "  +1.12%  "
$ws.Range("D12").Value = "'4.633"
$ws.Range("E12").Value = "  +0.25%  "
$ws.Range("D13").Value = "1.707.74"
$ws.Range("E13").Value = "  +1.38%  "
$ws.Range("D14").Value = "1.959.62"
$ws.Range("E14").Value = "  +4.01%  "
$ws.Range("D15").Value = "'0.5867"
$ws.Range("E15").Value = "  +4.54%  "
$ws.Range("D16").Value = "0.0₅8255"
$ws.Range("E16").Value = "  +0.98%  "
$ws.Range("D17").Value = "'67.92"
$ws.Range("E17").Value = "  +3.63%  "
$ws.Range("D18").Value = "27.483.37"
$ws.Range("E18").Value = "  +5.09%  "
$ws.Range("D19").Value = "'222.81"
$ws.Range("E19").Value = "  +15.63%  "
$ws.Range("D20").Value = "'1.005"
$ws.Range("E20").Value = "  +0.12%  "
$ws.Range("D21").Value = "'4.734"
$ws.Range("E21").Value = "  +1.75%  "
$ws.Range("D22").Value = "'10.68"
$ws.Range("E22").Value = "  +1.93%  "
$ws.Range("D23").Value = "'6.097"
$ws.Range("E23").Value = "  +2.35%  "
$ws.Range("D24").Value = "'1.006"
$ws.Range("E24").Value = "  +0.20%  "
$ws.Range("D25").Value = "'147.96"
$ws.Range("E25").Value = "  +1.82%  "
$ws.Range("D26").Value = "'0.1231"
$ws.Range("E26").Value = "  +3.05%  "
$ws.Range("E27").Value = "  +11.05%  "
$ws.Range("D28").Value = "'7.405"
$ws.Range("E28").Value = "  +1.94%  "
$ws.Range("E29").Value = "  +4.62%  "
$ws.Range("D30").Value = "'0.05554"
$ws.Range("E30").Value = "  +1.35%  "
$ws.Range("E31").Value = "  +2.55%  "
$ws.Range("D32").Value = "'3.568"
$ws.Range("E32").Value = "  +2.83%  "
$ws.Range("D33").Value = "'3.459"
$ws.Range("E33").Value = "  +2.63%  "
$ws.Range("D34").Value = "'1.658"
$ws.Range("E34").Value = "  +6.09%  "
$ws.Range("D35").Value = "'2.452"
$ws.Range("E35").Value = "  +2.08%  "
$ws.Range("D36").Value = "'0.9580"
$ws.Range("E36").Value = "  +0.88%  "
$ws.Range("D37").Value = "'2.817"
$ws.Range("E37").Value = "  +1.26%  "
$ws.Range("D38").Value = "'0.5917"
$ws.Range("E38").Value = "  +4.06%  "
$ws.Range("D39").Value = "'0.01642"
$ws.Range("E39").Value = "  +3.22%  "
$ws.Range("D40").Value = "'5.866"
$ws.Range("E40").Value = "  -0.04%  "
$ws.Range("D41").Value = "'0.8557"
$ws.Range("E41").Value = "  +3.01%  "
$ws.Range("D42").Value = "1.055.16"
$ws.Range("E42").Value = "  +2.39%  "
$ws.Range("D43").Value = "'1.005"
$ws.Range("E43").Value = "  +0.16%  "
$ws.Range("D44").Value = "'101.50"
$ws.Range("E44").Value = "  +0.53%  "
$ws.Range("D45").Value = "1.865.50"
$ws.Range("E45").Value = "  +3.94%  "
$ws.Range("D46").Value = "0.0₈113"
$ws.Range("E46").Value = "  +8.79%  "
$ws.Range("D47").Value = "'58.91"
$ws.Range("E47").Value = "  +1.74%  "
$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D48").Value = "'0.4446"
$ws.Range("E48").Value = "  +2.28%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "'8.190"
$ws.Range("E49").Value = "  +2.22%  "
$ws.Range("D50").Value = "'1.000"
$ws.Range("E50").Value = "  +0.32%  "
$ws.Range("D51").Value = "'0.05279"
$ws.Range("E51").Value = "  +1.26%  "

foreach ($addr in @("D5","D6","D7","D8","D9","D10","D11","D12","D15","D17","D19","D20","D21","D22","D23","D24","D25","D26","D28","D30","D32","D33","D34","D35","D36","D37","D38","D39","D40","D41","D43","D44","D47","D48","D49","D50","D51")) {
    $ws.Range($addr).Style = "Normal"
}

Write-Host "Applied 99 cell updates"
